$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44186
$ws.Range("M2").Value = 40
$ws.Range("N2").Value = 15000
$ws.Range("O2").Value = 15000
$ws.Range("P2").Value = 15000
$ws.Range("S2").Value = 3000

# Row 3
$ws.Range("D3").Value = 44179
$ws.Range("M3").Value = 45
$ws.Range("N3").Value = 20000
$ws.Range("O3").Value = 20000
$ws.Range("P3").Value = 20000
$ws.Range("S3").Value = 4000

# Row 4
$ws.Range("D4").Value = 44196
$ws.Range("M4").Value = 56
$ws.Range("N4").Value = 15000
$ws.Range("O4").Value = 15000
$ws.Range("P4").Value = 15000
$ws.Range("S4").Value = 3000

# Row 7
$ws.Range("D7").Value = 44902
$ws.Range("M7").Value = 35
$ws.Range("N7").Value = 12000
$ws.Range("O7").Value = 12000
$ws.Range("P7").Value = 12000
$ws.Range("S7").Value = 2400

# Row 8
$ws.Range("D8").Value = 44931
$ws.Range("M8").Value = 50
$ws.Range("N8").Value = 18000
$ws.Range("O8").Value = 18000
$ws.Range("P8").Value = 18000
$ws.Range("S8").Value = 3600

# Row 9
$ws.Range("D9").Value = 44175
$ws.Range("M9").Value = 25

# Row 10
$ws.Range("D10").Value = 44914
$ws.Range("M10").Value = 56
$ws.Range("N10").Value = 23000
$ws.Range("O10").Value = 23000
$ws.Range("P10").Value = 23000
$ws.Range("S10").Value = 4600

# Row 11
$ws.Range("D11").Value = 44181
$ws.Range("M11").Value = 30
$ws.Range("N11").Value = 20000
$ws.Range("O11").Value = 20000
$ws.Range("P11").Value = 20000
$ws.Range("S11").Value = 4000

# Row 12
$ws.Range("D12").Value = 44193
$ws.Range("M12").Value = 40
$ws.Range("N12").Value = 15000
$ws.Range("O12").Value = 15000
$ws.Range("P12").Value = 15000
$ws.Range("S12").Value = 3000

# Row 13
$ws.Range("D13").Value = 44188
$ws.Range("M13").Value = 30

